# ------------------------------------------------------------------
# Edit script: Image processing visualization system.pptx
#
# Summary of changes applied (per authoritative diff):
#  1. Delete the two trailing slides ("3.2 Development timeline" and
#     "3.3 Tasks and milestones", originally slides 11 and 12).
#  2. Update the auto-generated "today" date field text that appears
#     on the slide-layout / notes-master date placeholders from
#     2023/4/12 to 2023/3/30.
#  3. Rename "3.1 Development schedule" -> "3 Development schedule"
#     on both slides that carry that heading.
#  4. On the "Requirements overview" slide, bold the short title
#     portion of each of the three numbered requirement bullets and
#     grow the containing text box to the new height.
#  5. On the "Membership arrangement" slide, change two occurrences
#     of "Project manager" to "Development engineer" (each preceded
#     by a newly inserted, specially coloured space run) and grow the
#     two containing text boxes to their new heights.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --------------------------------------------------------------
# 1. Remove the last two slides (SlideID 306 and 307).
# --------------------------------------------------------------
$totalSlides = $p.Slides.Count
$p.Slides.Item($totalSlides).Delete()
$p.Slides.Item($totalSlides - 1).Delete()

# --------------------------------------------------------------
# 2. Update the date placeholder text (2023/4/12 -> 2023/3/30)
#    across every slide layout and the notes master.
# --------------------------------------------------------------
$oldDate = "2023/4/12"
$newDate = "2023/3/30"

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($sj = 1; $sj -le $layout.Shapes.Count; $sj++) {
        $lsh = $layout.Shapes.Item($sj)
        if ($lsh.HasTextFrame) {
            if ($lsh.TextFrame.TextRange.Text -eq $oldDate) {
                $lsh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$notesMaster = $p.NotesMaster
for ($sj = 1; $sj -le $notesMaster.Shapes.Count; $sj++) {
    $nsh = $notesMaster.Shapes.Item($sj)
    if ($nsh.HasTextFrame) {
        if ($nsh.TextFrame.TextRange.Text -eq $oldDate) {
            $nsh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --------------------------------------------------------------
# Helper: find the shape on a slide whose text frame text matches
# exactly, returning its index (0 if not found).
# --------------------------------------------------------------
function Find-ShapeByText($slide, $text) {
    for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
        $cand = $slide.Shapes.Item($k)
        if ($cand.HasTextFrame) {
            if ($cand.TextFrame.TextRange.Text -eq $text) {
                return $k
            }
        }
    }
    return 0
}

# --------------------------------------------------------------
# 3. "3.1 Development schedule" -> "3 Development schedule"
#    (appears on the development-schedule slides).
# --------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $idx = Find-ShapeByText $slide "3.1 Development schedule"
    if ($idx -gt 0) {
        $titleShape = $slide.Shapes.Item($idx)
        $origHeight = $titleShape.Height
        $titleShape.TextFrame.TextRange.Text = "3 Development schedule"
        # The shape auto-fits its text; restore the original height so
        # the box geometry is left untouched, matching the source edit.
        $titleShape.Height = $origHeight
    }
}

# --------------------------------------------------------------
# 4. Requirements overview slide: bold the bullet "titles".
# --------------------------------------------------------------
$reqSlideIdx = 0
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $idx = Find-ShapeByText $slide "1.2 Requirements overview"
    if ($idx -gt 0) {
        $reqSlideIdx = $si
    }
}

if ($reqSlideIdx -gt 0) {
    $reqSlide = $p.Slides.Item($reqSlideIdx)
    $bodyIdx = Find-ShapeByText $reqSlide "1.2 Requirements overview"
    # the body textbox is the other big shape holding all 3 bullets
    for ($k = 1; $k -le $reqSlide.Shapes.Count; $k++) {
        $cand = $reqSlide.Shapes.Item($k)
        if ($cand.HasTextFrame) {
            if ($cand.TextFrame.TextRange.Text.Length -gt 200) {
                $bodyShapeIdx = $k
            }
        }
    }

    $bodyShape = $reqSlide.Shapes.Item($bodyShapeIdx)
    $tr = $bodyShape.TextFrame.TextRange

    $titles = @(
        @{ Prefix = "1. "; Title = "Visualization of image enhancement method" },
        @{ Prefix = "2. "; Title = "Visualization of convolutional network" },
        @{ Prefix = "3. "; Title = "Simple image recognition display" }
    )

    $paraCount = $tr.Paragraphs().Count
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $para = $tr.Paragraphs($pi)
        foreach ($item in $titles) {
            $full = $item.Prefix + $item.Title
            if ($para.Text.StartsWith($full)) {
                $sub = $para.Characters(($item.Prefix.Length + 1), $item.Title.Length)
                $sub.Font.Bold = $true
            }
        }
    }

    # Grow the text box to the new height (EMU 4112895 -> 4397358).
    $bodyShape.Height = 346.2486724853516
}

# --------------------------------------------------------------
# 5. Membership arrangement slide: "Project manager" ->
#    "Development engineer" for the two week4-10 entries.
# --------------------------------------------------------------
$memSlideIdx = 0
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $idx = Find-ShapeByText $slide "1.3 Membership arrangement"
    if ($idx -gt 0) {
        $memSlideIdx = $si
    }
}

if ($memSlideIdx -gt 0) {
    $memSlide = $p.Slides.Item($memSlideIdx)

    # Locate the two target text boxes precisely by position, since
    # several boxes share identical text ("Project manager").
    $boxAIdx = 0
    $boxBIdx = 0
    for ($k = 1; $k -le $memSlide.Shapes.Count; $k++) {
        $cand = $memSlide.Shapes.Item($k)
        if ($cand.HasTextFrame) {
            $leftEmu = [Math]::Round($cand.Left * 12700)
            $topEmu = [Math]::Round($cand.Top * 12700)
            if ($leftEmu -eq 0 -and $topEmu -eq 3350895) {
                $boxAIdx = $k
            }
            if ($leftEmu -eq 6127115 -and $topEmu -eq 3319780) {
                $boxBIdx = $k
            }
        }
    }

    # --- Box A (top-left block): week4-10 "Project manager" -> "Development engineer"
    $boxA = $memSlide.Shapes.Item($boxAIdx)
    $trA = $boxA.TextFrame.TextRange
    for ($pi = 1; $pi -le $trA.Paragraphs().Count; $pi++) {
        $para = $trA.Paragraphs($pi)
        if ($para.Text.StartsWith("week4-10:") -and $para.Text.Contains("Project manager")) {
            $pmStart = $para.Text.IndexOf("Project manager") + 1
            $pmRange = $para.Characters($pmStart, "Project manager".Length)
            $pmRange.InsertBefore(" ")

            $pmStart2 = $para.Text.IndexOf("Project manager") + 1
            $spaceRange = $para.Characters(($pmStart2 - 1), 1)
            $spaceRange.Font.Color.RGB = 0xC17015

            $pmStart3 = $para.Text.IndexOf("Project manager") + 1
            $pmRange3 = $para.Characters($pmStart3, "Project manager".Length)
            $pmRange3.Text = "Development engineer"
        }
    }
    $boxA.Height = 180.6824417114258

    # --- Box B (week7-9 block): second week4-10 "Project manager" -> "Development engineer "
    $boxB = $memSlide.Shapes.Item($boxBIdx)
    $trB = $boxB.TextFrame.TextRange
    for ($pi = 1; $pi -le $trB.Paragraphs().Count; $pi++) {
        $para = $trB.Paragraphs($pi)
        if ($para.Text.StartsWith("week4-10:") -and $para.Text.Contains("Project manager")) {
            $pmStart = $para.Text.IndexOf("Project manager") + 1
            $pmRange = $para.Characters($pmStart, "Project manager".Length)
            $pmRange.InsertBefore(" ")

            $pmStart2 = $para.Text.IndexOf("Project manager") + 1
            $spaceRange = $para.Characters(($pmStart2 - 1), 1)
            $spaceRange.Font.Color.RGB = 0xC17015

            $pmStart3 = $para.Text.IndexOf("Project manager") + 1
            $pmRange3 = $para.Characters($pmStart3, "Project manager".Length)
            $pmRange3.Text = "Development engineer "
        }
    }
    $boxB.Height = 268.81984251968503
}
